# Update countries & provincias Spain
# - Ucrania overtakes Polonia in the ranking (row 36 becomes Ucrania, row 37 becomes Polonia)
# - Refresh several countries' case counters (Israel, Lituania, Eslovaquia, Georgia)
# - Bump the "last updated" timestamp from 09:05 to 09:35

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: now Ucrania (was Polonia), with updated figures
$ws.Cells.Item(36, 1).Value = "Ucrania"
$ws.Cells.Item(36, 2).Value = 20986
$ws.Cells.Item(36, 3).Value = 406
$ws.Cells.Item(36, 4).Value = 7108
$ws.Cells.Item(36, 5).Value = 13261
$ws.Cells.Item(36, 7).Value = 12
$ws.Cells.Item(36, 8).Value = 617

# Row 37: now Polonia (was Ucrania), with updated figures
$ws.Cells.Item(37, 1).Value = "Polonia"
$ws.Cells.Item(37, 2).Value = 20931
$ws.Cells.Item(37, 4).Value = 9194
$ws.Cells.Item(37, 5).Value = 10744
$ws.Cells.Item(37, 8).Value = 993

# Row 41: Israel - Recuperados/Casos activos updated
$ws.Cells.Item(41, 4).Value = 14093
$ws.Cells.Item(41, 5).Value = 2340

# Row 94: Lituania - updated figures
$ws.Cells.Item(94, 2).Value = 1623
$ws.Cells.Item(94, 3).Value = 7
$ws.Cells.Item(94, 4).Value = 1138
$ws.Cells.Item(94, 5).Value = 422

# Row 97: Eslovaquia - updated figures
$ws.Cells.Item(97, 2).Value = 1509
$ws.Cells.Item(97, 3).Value = 5
$ws.Cells.Item(97, 4).Value = 1301
$ws.Cells.Item(97, 5).Value = 180

# Row 122: Georgia - updated figures
$ws.Cells.Item(122, 2).Value = 730
$ws.Cells.Item(122, 3).Value = 2
$ws.Cells.Item(122, 4).Value = 522
$ws.Cells.Item(122, 5).Value = 196

# Update the "last refreshed" timestamp shown in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Mayo de 2020 a las 09:35"
